$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 48 ("Primera"/"Segunda" pair),
# pushing all existing data rows (old 48-120) down to 50-122.
$ws.Rows("48:49").Insert()

# New row 48: Primera, fecha 45210 (2023-10-11), volumen 40, precios 16000, precio/kg 1600
$row48 = @(3, "Femacal de La Calera", "Coquimbo", 45210, 5, "Fruta", 100108, "Tropicales y subtropicales", 100108004, "Papaya", "Cultivar IV Región", "Primera", 40, 16000, 16000, 16000, "`$/bandeja 10 kilos", "Provincia del Elquí", 1600, 10)
for ($i = 0; $i -lt $row48.Length; $i++) {
    $ws.Cells.Item(48, $i + 1).Value = $row48[$i]
}

# New row 49: Segunda, fecha 45210 (2023-10-11), volumen 48, precios 13000, precio/kg 1300
$row49 = @(3, "Femacal de La Calera", "Coquimbo", 45210, 5, "Fruta", 100108, "Tropicales y subtropicales", 100108004, "Papaya", "Cultivar IV Región", "Segunda", 48, 13000, 13000, 13000, "`$/bandeja 10 kilos", "Provincia del Elquí", 1300, 10)
for ($i = 0; $i -lt $row49.Length; $i++) {
    $ws.Cells.Item(49, $i + 1).Value = $row49[$i]
}
